$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Put" request row, mirroring the existing Get/Post rows:
# A4 = "Put" (request name), B4 = endpoint URL, C4 = status code 200.
$ws.Range("A4").Value = "Put"
$ws.Range("B4").Value = "https://bookstore.toolsqa.com"
$ws.Range("C4").Value = 200

# Copy formatting from row 2 (A2/C2 plain bordered style, B2 left/vcenter style)
# so the new row matches the sheet's existing look.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

$ws.Range("C6").Select() | Out-Null
